$d = $word.ActiveDocument

# 1) "Objectives: Make a user-friendly GUI" -> "Objectives:"
$d.Content.Find.Execute("Objectives: Make a user-friendly GUI", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Objectives:", 2)

# 2) "Extract data from a csv file" -> "Research similar projects" (the trailing "." run is untouched)
$d.Content.Find.Execute("Extract data from a csv file", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Research similar projects", 2)

# 3) "Make tables from data showing each individual patient." -> "Roles Assignment" + " " + " "
$p6 = $d.Paragraphs.Item(6)
$p6.Range.Text = "Roles Assignment"
$p6 = $d.Paragraphs.Item(6)
$p6.Range.InsertAfter(" ")
$p6 = $d.Paragraphs.Item(6)
$p6.Range.InsertAfter(" ")

# 4) "Generate reports on individual patients." -> "Design" + " " + "(Figure out technical requirements and evaluate the use cases)"
$p7 = $d.Paragraphs.Item(7)
$p7.Range.Text = "Design"
$p7 = $d.Paragraphs.Item(7)
$p7.Range.InsertAfter(" ")
$p7 = $d.Paragraphs.Item(7)
$p7.Range.InsertAfter("(Figure out technical requirements and evaluate the use cases)")

# 5) Replace the trailing (empty) paragraph with four new paragraphs. Reuse the
#    existing trailing empty paragraph as the first of the four (so no extra
#    paragraph is left over), then create the other three (still empty) so
#    paragraph-level formatting (tab stops) added afterwards does not get
#    inherited by the next InsertParagraphAfter.
$p8 = $d.Paragraphs.Item(8)          # paragraph 8 (was the trailing empty paragraph): "Test planning."
$p8.Range.InsertParagraphAfter()     # paragraph 9: "Implementation and Testing"
$p9 = $d.Paragraphs.Item(9)
$p9.Range.InsertParagraphAfter()     # paragraph 10: "Write the report."
$p10 = $d.Paragraphs.Item(10)
$p10.Range.InsertParagraphAfter()    # paragraph 11: tab only

# paragraph 8: "Test " + "planning."
$p8 = $d.Paragraphs.Item(8)
$p8.Range.Text = "Test "
$p8 = $d.Paragraphs.Item(8)
$p8.Range.InsertAfter("planning.")

# paragraph 9: "Implementation" + " and Testing", with a centered tab stop at 4513 twips (225.65 pt)
$p9 = $d.Paragraphs.Item(9)
$p9.Range.Text = "Implementation"
$p9 = $d.Paragraphs.Item(9)
$p9.Range.InsertAfter(" and Testing")
$p9 = $d.Paragraphs.Item(9)
$p9.Range.ParagraphFormat.TabStops.Add(225.65, 1)

# paragraph 10: "Write the report.", with the same centered tab stop
$p10 = $d.Paragraphs.Item(10)
$p10.Range.Text = "Write the report."
$p10 = $d.Paragraphs.Item(10)
$p10.Range.ParagraphFormat.TabStops.Add(225.65, 1)

# paragraph 11: a single tab character, with the same centered tab stop
$p11 = $d.Paragraphs.Item(11)
$p11.Range.ParagraphFormat.TabStops.Add(225.65, 1)
$p11 = $d.Paragraphs.Item(11)
$p11.Range.Text = [char]9

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    Write-Output "$i -> [$($pp.Range.Text)]"
}
